$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '58.088.43'
$ws.Cells.Item(2, 5).Value = '  +2.87%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.354.23'
$ws.Cells.Item(3, 5).Value = '  +1.31%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.05%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '545.21'
$ws.Cells.Item(5, 5).Value = '  +6.49%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '134.64'
$ws.Cells.Item(6, 5).Value = '  +2.52%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.999'
$ws.Cells.Item(7, 5).Value = '  +0.04%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.539'
$ws.Cells.Item(8, 5).Value = '  +1.07%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +1.18%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.90%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '5.42'
$ws.Cells.Item(12, 5).Value = '  +3.30%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.360'
$ws.Cells.Item(13, 5).Value = '  +7.24%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2.771.78'
$ws.Cells.Item(14, 5).Value = '  +1.22%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '23.59'
$ws.Cells.Item(15, 5).Value = '  +0.50%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '58.057.19'
$ws.Cells.Item(16, 5).Value = '  +2.83%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +1.26%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.352.43'
$ws.Cells.Item(18, 5).Value = '  +1.47%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '10.59'
$ws.Cells.Item(19, 5).Value = '  +2.03%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '333.91'
$ws.Cells.Item(20, 5).Value = '  +1.95%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +2.13%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.74'
$ws.Cells.Item(22, 5).Value = '  +0.42%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.32%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '61.87'
$ws.Cells.Item(24, 5).Value = '  +0.98%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +4.39%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.48'
$ws.Cells.Item(26, 5).Value = '  -1.05%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.996'
$ws.Cells.Item(27, 5).Value = '  -0.30%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.40'
$ws.Cells.Item(28, 5).Value = '  +7.92%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +5.68%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '169.83'
$ws.Cells.Item(30, 5).Value = '  +1.23%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.0₃0734'
$ws.Cells.Item(31, 5).Value = '  +2.38%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '6.16'
$ws.Cells.Item(32, 5).Value = '  +1.03%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.04'
$ws.Cells.Item(33, 5).Value = '  +17.79%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '18.50'
$ws.Cells.Item(34, 5).Value = '  +0.97%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  +0.04%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.999'
$ws.Cells.Item(36, 5).Value = '  +0.10%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'NEARProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.20'
$ws.Cells.Item(37, 5).Value = '  +7.19%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'ImmutableX'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.27'
$ws.Cells.Item(38, 5).Value = '  +1.58%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.65'
$ws.Cells.Item(39, 5).Value = '  +6.00%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '39.23'
$ws.Cells.Item(40, 5).Value = '  +1.64%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '149.59'
$ws.Cells.Item(41, 5).Value = '  +0.78%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.381'
$ws.Cells.Item(42, 5).Value = '  +2.06%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '288.41'
$ws.Cells.Item(43, 5).Value = '  +4.52%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +1.67%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '19.32'
$ws.Cells.Item(45, 5).Value = '  +6.45%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0927'
$ws.Cells.Item(46, 5).Value = '  -0.02%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0506'
$ws.Cells.Item(47, 5).Value = '  +2.54%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.561'
$ws.Cells.Item(48, 5).Value = '  +0.97%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +1.73%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '17.58'
$ws.Cells.Item(50, 5).Value = '  +3.74%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Polygon'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.381'
$ws.Cells.Item(51, 5).Value = '  +0.20%  '

